# Auto-generated edit script applying the Seraph_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 300
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H51").Value = 79234.57000000001
$ws.Range("I51").Value = 11281.5
$ws.Range("J51").Value = 130199.375
$ws.Range("K51").Value = 11281.5
$ws.Range("L51").Value = 130199.375
$ws.Range("M51").Value = -10797.5
$ws.Range("N51").Value = -131167.375
$ws.Range("H58").Value = 6615
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 6988.5713
$ws.Range("K58").Value = 12000
$ws.Range("L58").Value = 20965.7139
$ws.Range("N58").Value = -21265.7139
$ws.Range("M58").Value = -11850
$ws.Range("H70").Value = 52054.785
$ws.Range("I70").Value = 4874.5
$ws.Range("J70").Value = 70926.89999999999
$ws.Range("K70").Value = 14623.5
$ws.Range("L70").Value = 212780.7
$ws.Range("M70").Value = -14353.5
$ws.Range("N70").Value = -213320.7
$ws.Range("H73").Value = 52054.785
$ws.Range("I73").Value = 4874.5
$ws.Range("J73").Value = 70926.89999999999
$ws.Range("K73").Value = 14623.5
$ws.Range("L73").Value = 212780.7
$ws.Range("M73").Value = -13687.5
$ws.Range("N73").Value = -214652.7
$ws.Range("H138").Value = 1877.174
$ws.Range("I138").Value = 1708.75
$ws.Range("K138").Value = 5126.25
$ws.Range("M138").Value = 13.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H61").Value = 13390.2
$ws.Range("I61").Value = 13390.2
$ws.Range("K61").Value = 13390.2
$ws.Range("M61").Value = -13178.2
$ws.Range("H112").Value = 57749.5
$ws.Range("J112").Value = 57749.5
$ws.Range("L112").Value = 57749.5
$ws.Range("N112").Value = -60703.5
$ws.Range("H132").Value = 3960.6155
$ws.Range("I132").Value = 4571
$ws.Range("K132").Value = 13713
$ws.Range("M132").Value = -11183
$ws.Range("H136").Value = 13390.2
$ws.Range("I136").Value = 13390.2
$ws.Range("K136").Value = 40170.60000000001
$ws.Range("M136").Value = -37620.60000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 16668658
$ws.Range("I7").Value = 25000488
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 25000488
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -25000375
$ws.Range("N7").Value = -5226
$ws.Range("H86").Value = 3121.75
$ws.Range("I86").Value = 2244.5
$ws.Range("K86").Value = 2244.5
$ws.Range("M86").Value = -1121.5
$ws.Range("H89").Value = 3121.75
$ws.Range("I89").Value = 2244.5
$ws.Range("K89").Value = 11222.5
$ws.Range("M89").Value = -5606.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 165.26666
$ws.Range("I7").Value = 140.8
$ws.Range("K7").Value = 140.8
$ws.Range("M7").Value = -27.80000000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35037184
$ws.Range("I4").Value = 43122490
$ws.Range("J4").Value = 881.5
$ws.Range("K4").Value = 129367470
$ws.Range("L4").Value = 2644.5
$ws.Range("M4").Value = -129367358
$ws.Range("N4").Value = -2868.5
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1327
$ws.Range("N16").Value = -3346
$ws.Range("H40").Value = 130.22223
$ws.Range("I40").Value = 159
$ws.Range("K40").Value = 636
$ws.Range("M40").Value = -567
$ws.Range("H86").Value = 275
$ws.Range("J86").Value = 275
$ws.Range("L86").Value = 825
$ws.Range("N86").Value = -3197
$ws.Range("H89").Value = 275
$ws.Range("J89").Value = 275
$ws.Range("L89").Value = 2475
$ws.Range("N89").Value = -14331
$ws.Range("H132").Value = 3751.4614
$ws.Range("J132").Value = 3666.3333
$ws.Range("L132").Value = 32996.9997
$ws.Range("N132").Value = -38056.9997
$ws.Range("H133").Value = 3500
$ws.Range("I133").Value = 3500
$ws.Range("K133").Value = 10500
$ws.Range("M133").Value = -5440
$ws.Range("H137").Value = 4352.9
$ws.Range("J137").Value = 4975
$ws.Range("L137").Value = 14925
$ws.Range("N137").Value = -25125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9999
$ws.Range("I46").Value = 9999
$ws.Range("K46").Value = 9999
$ws.Range("M46").Value = -9843
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2495.4
$ws.Range("I46").Value = 2245.6667
$ws.Range("J46").Value = 2870
$ws.Range("K46").Value = 2245.6667
$ws.Range("L46").Value = 2870
$ws.Range("M46").Value = -2057.6667
$ws.Range("N46").Value = -3246
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 1638.8
$ws.Range("I93").Value = 898
$ws.Range("K93").Value = 898
$ws.Range("M93").Value = 350
$ws.Range("H136").Value = 6792.5
$ws.Range("I136").Value = 5856
$ws.Range("J136").Value = 8665.5
$ws.Range("K136").Value = 17568
$ws.Range("L136").Value = 25996.5
$ws.Range("M136").Value = -15018
$ws.Range("N136").Value = -31096.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 516875
$ws.Range("J26").Value = 516875
$ws.Range("L26").Value = 516875
$ws.Range("N26").Value = -517461
$ws.Range("H47").Value = 44665
$ws.Range("I47").Value = 44000
$ws.Range("K47").Value = 44000
$ws.Range("M47").Value = -43428
$ws.Range("H62").Value = 4685.4
$ws.Range("I62").Value = 3822.5
$ws.Range("J62").Value = 4999.1816
$ws.Range("K62").Value = 3822.5
$ws.Range("L62").Value = 4999.1816
$ws.Range("M62").Value = -3198.5
$ws.Range("N62").Value = -6247.1816
$ws.Range("H65").Value = 4685.4
$ws.Range("I65").Value = 3822.5
$ws.Range("J65").Value = 4999.1816
$ws.Range("K65").Value = 19112.5
$ws.Range("L65").Value = 24995.908
$ws.Range("M65").Value = -15992.5
$ws.Range("N65").Value = -31235.908
$ws.Range("H70").Value = 52500
$ws.Range("J70").Value = 52500
$ws.Range("L70").Value = 52500
$ws.Range("N70").Value = -53130
$ws.Range("H73").Value = 52500
$ws.Range("J73").Value = 52500
$ws.Range("L73").Value = 52500
$ws.Range("N73").Value = -54684
